# Initial commit of backend
# Adds a new "Other" worksheet (after "Items") with GitHub / CircleCI links,
# and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$itemsSheet = $wb.Worksheets.Item("Items")

# Add a new worksheet right after "Items".
$otherSheet = $wb.Worksheets.Add($null, $itemsSheet)
$otherSheet.Name = "Other"

# Populate in the same order the original author did, so new shared
# strings land in the same sequence (CircleCI, CircleCI URL, then GitHub).
$otherSheet.Range("A2").Value = "CircleCI"
$otherSheet.Range("B2").Value = "https://app.circleci.com/pipelines/github/jpickup/udacity-capstone"
$otherSheet.Range("A1").Value = "GitHub"
$otherSheet.Range("B1").Value = "https://github.com/jpickup/udacity-capstone"

# Make the new sheet the active/selected tab, with A3 selected.
$otherSheet.Activate()
$otherSheet.Range("A3").Select() | Out-Null
